$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data (and swap TheSandbox / InternetComputer rows 39-40)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.817.07"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.763.10"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.21"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4252"
$ws.Range("E7").Value = "  -4.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3633"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07527"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.52"
$ws.Range("E10").Value = "  -5.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.092"
$ws.Range("E11").Value = "  -3.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.75"
$ws.Range("E13").Value = "  -5.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.071"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.280"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.766.61"
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.37"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001064"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06384"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.05"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.911"
$ws.Range("E22").Value = "  -5.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.856.70"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.117"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.12"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.26"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.984.91"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.153"
$ws.Range("E29").Value = "  -7.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.14"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("E31").Value = "  -6.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.683"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08889"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.27"
$ws.Range("E35").Value = "  -6.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02290"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2104"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06035"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.977"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6342"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.175"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9993"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.909"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.398"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.31"
$ws.Range("E45").Value = "  -5.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5879"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.689"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.990"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.82"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.185"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06834"
$ws.Range("E51").Value = "  -2.28%  "
